$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text cell (matches the source inline-string
    # data) instead of letting Excel auto-coerce numeric- or date-looking
    # strings into Number/Date cells, then restore the original "Normal"
    # style so no stray number-format style sticks to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# M3: cited_by_count 2 -> 3
Set-TextValue $ws.Range("M3") "3"

# Rows 6 and 7 got swapped (their DOM/Banner author ids were reassigned to
# the correct rows). Columns J, K, L, M, N, O, Q hold identical values in
# both rows already, so only A-I and P need to be swapped.
$row6 = @($ws.Range("A6").Text, $ws.Range("B6").Text, $ws.Range("C6").Text, $ws.Range("D6").Text, $ws.Range("E6").Text, $ws.Range("F6").Text, $ws.Range("G6").Text, $ws.Range("H6").Text, $ws.Range("I6").Text, $ws.Range("P6").Text)
$row7 = @($ws.Range("A7").Text, $ws.Range("B7").Text, $ws.Range("C7").Text, $ws.Range("D7").Text, $ws.Range("E7").Text, $ws.Range("F7").Text, $ws.Range("G7").Text, $ws.Range("H7").Text, $ws.Range("I7").Text, $ws.Range("P7").Text)

Set-TextValue $ws.Range("A6") $row7[0]
Set-TextValue $ws.Range("B6") $row7[1]
Set-TextValue $ws.Range("C6") $row7[2]
Set-TextValue $ws.Range("D6") $row7[3]
Set-TextValue $ws.Range("E6") $row7[4]
Set-TextValue $ws.Range("F6") $row7[5]
Set-TextValue $ws.Range("G6") $row7[6]
Set-TextValue $ws.Range("H6") $row7[7]
Set-TextValue $ws.Range("I6") $row7[8]
Set-TextValue $ws.Range("P6") $row7[9]

Set-TextValue $ws.Range("A7") $row6[0]
Set-TextValue $ws.Range("B7") $row6[1]
Set-TextValue $ws.Range("C7") $row6[2]
Set-TextValue $ws.Range("D7") $row6[3]
Set-TextValue $ws.Range("E7") $row6[4]
Set-TextValue $ws.Range("F7") $row6[5]
Set-TextValue $ws.Range("G7") $row6[6]
Set-TextValue $ws.Range("H7") $row6[7]
Set-TextValue $ws.Range("I7") $row6[8]
Set-TextValue $ws.Range("P7") $row6[9]
